$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new report rows (4, 5, 6) below the existing data, matching
# the existing table's layout (Report Id, Popular Item, Popular Customer,
# Popular Staff, Created Date). The "Created Date" column keeps the same
# custom date format used by the rows above it.

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 45410
$ws.Range("E4").NumberFormat = "yyyy-MM-dd"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 45410
$ws.Range("E5").NumberFormat = "yyyy-MM-dd"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 45410
$ws.Range("E6").NumberFormat = "yyyy-MM-dd"

Write-Output "added rows 4-6"
